$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as text so exact formatting
# (leading/trailing zeros, thousand-dot separators, etc.) is preserved.
$textCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D11', 'D13', 'D15', 'D16', 'D17', 'D18', 'D20', 'D23', 'D25', 'D27', 'D28', 'D29', 'D30', 'D32', 'D33', 'D35', 'D37', 'D39', 'D42', 'D43', 'D44', 'D48', 'D50')
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.404.61'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').Value = '3.391.37'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '587.46'
$ws.Range('E5').Value = '  +1.02%  '
$ws.Range('D6').Value = '179.84'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('E9').Value = '  +6.10%  '
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('D11').Value = '48.50'
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('E12').Value = '  +3.20%  '
$ws.Range('D13').Value = '678.48'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D15').Value = '3.937.82'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = '69.455.25'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.392.72'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.120'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').Value = '11.27'
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').Value = '17.10'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  +3.82%  '
$ws.Range('D25').Value = '3.93'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('E26').Value = '  +1.27%  '
$ws.Range('D27').Value = '9.68'
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('D28').Value = '34.14'
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('D29').Value = '8.71'
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('D32').Value = '556.19'
$ws.Range('E32').Value = '  -1.71%  '
$ws.Range('D33').Value = '3.62'
$ws.Range('E33').Value = '  +6.50%  '
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('D35').Value = '58.03'
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = '3.684.77'
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('E38').Value = '  +5.42%  '
$ws.Range('D39').Value = '35.02'
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').Value = '0.0₃0699'
$ws.Range('E42').Value = '  +3.48%  '
$ws.Range('D43').Value = '0.339'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').Value = '0.0423'
$ws.Range('E44').Value = '  +4.18%  '
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('D48').Value = '1.42'
$ws.Range('E48').Value = '  +6.14%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').Value = '132.63'
$ws.Range('E51').Value = '  +3.28%  '
